$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Status'
$ws.Range("B1").Value = 'docType'
$ws.Range("C1").Value = 'id'
$ws.Range("D1").Value = 'title'
$ws.Range("E1").Value = 'description'
$ws.Range("F1").Value = 'slug'
$ws.Range("G1").Value = 'date'
$ws.Range("H1").Value = 'author'
$ws.Range("I1").Value = 'image'
$ws.Range("J1").Value = 'category'
$ws.Range("K1").Value = 'tags'
$ws.Range("L1").Value = 'relatedDoc1'
$ws.Range("M1").Value = 'relatedDoc2'
$ws.Range("N1").Value = 'relatedDoc3'
$ws.Range("O1").Value = 'relatedDoc4'
$ws.Range("P1").Value = 'relatedDoc5'
$ws.Range("Q1").Value = 'editor_img'
$ws.Range("R1").Value = 'editor_bio'
$ws.Range("S1").Value = 'coverImage'

$ws.Range("A2").Value = '''''Pending'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = '''0'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = '''What is Portflorio?'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''Portflorio is a blog template built for Next.js. This is a free template that uses TypeScript and article management using Markdown, while maintaining the speed-up features that are the hallmark of Next.js.'
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = '''sample.md'
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = '''2023-12-31'
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = '''Shoto Morisaki'
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Value = '''sample/sample1.jpg'
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = '''Other'
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").Value = '''[Sample Markdown]'
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").Value = '''""'
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value = '''""'
$ws.Range("M2").Style = "Normal"
$ws.Range("N2").Value = '''""'
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").Value = '''""'
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").Value = '''""'
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Value = ''''''''
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").Value = ''''''''
$ws.Range("R2").Style = "Normal"
$ws.Range("S2").Value = ''''
$ws.Range("S2").Style = "Normal"

$ws.Range("A3").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = '''11'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = '''ACM (Association for Computing Machinery) Research Lab'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''This project revolves around implementing and understanding the REDCODER model, which is a code generation model. The team aims to evaluate the effectiveness of REDCODER for specific tasks by comparing its performance under different scenarios'
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = '''ACM_Research_Lab.md'
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = '''2024-4-1'
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = '''Shoto Morisaki'
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = '''project/6/ACM_logo.png'
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = '''Research'
$ws.Range("J3").Style = "Normal"
$ws.Range("K3").Value = '''[ LLM, Natural Language, LangChain, RAG ]'
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = '''""'
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").Value = '''""'
$ws.Range("M3").Style = "Normal"
$ws.Range("N3").Value = '''""'
$ws.Range("N3").Style = "Normal"
$ws.Range("O3").Value = '''""'
$ws.Range("O3").Style = "Normal"
$ws.Range("P3").Value = '''""'
$ws.Range("P3").Style = "Normal"
$ws.Range("Q3").Value = ''''''''
$ws.Range("Q3").Style = "Normal"
$ws.Range("R3").Value = ''''''''
$ws.Range("R3").Style = "Normal"
$ws.Range("S3").Value = ''''
$ws.Range("S3").Style = "Normal"

$ws.Range("A4").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = '''10'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = '''Software Engineering Internship at Tsubasa'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''Internship experience is a significant contributor to my development of engineering skills. Most of the work required me to use the serverless framework, and each time, there was a self-education aspect that I can leverage later in my career.'
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = '''SWE_Internship.md'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = '''2024-2-27'
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = '''Shoto Morisaki'
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = '''project/5/running.jpg'
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = '''Project'
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").Value = '''[ AWS, Docker, Python Flask, Material-UI, ApexCharts, TypeScript, React]'
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = '''""'
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Value = '''""'
$ws.Range("M4").Style = "Normal"
$ws.Range("N4").Value = '''""'
$ws.Range("N4").Style = "Normal"
$ws.Range("O4").Value = '''""'
$ws.Range("O4").Style = "Normal"
$ws.Range("P4").Value = '''""'
$ws.Range("P4").Style = "Normal"
$ws.Range("Q4").Value = ''''''''
$ws.Range("Q4").Style = "Normal"
$ws.Range("R4").Value = ''''''''
$ws.Range("R4").Style = "Normal"
$ws.Range("S4").Value = ''''
$ws.Range("S4").Style = "Normal"

$ws.Range("A5").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = '''15'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = '''Takeaway of 8 months Internship as Software Engineer'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''Wrapping up 8 months of internship at start up software engineering company as a Software Engineer has been an enriching journey. Here''s a summary of key takeaways'
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = '''Takeaway_from_SWE_internship.md'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = '''2024-02-20'
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = '''Shoto Morisaki'
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Value = '''sample/sample7.jpg'
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = '''Blog'
$ws.Range("J5").Style = "Normal"
$ws.Range("K5").Value = '''[ Internship, Software Engineering ]'
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = '''""'
$ws.Range("L5").Style = "Normal"
$ws.Range("M5").Value = '''""'
$ws.Range("M5").Style = "Normal"
$ws.Range("N5").Value = '''""'
$ws.Range("N5").Style = "Normal"
$ws.Range("O5").Value = '''""'
$ws.Range("O5").Style = "Normal"
$ws.Range("P5").Value = '''""'
$ws.Range("P5").Style = "Normal"
$ws.Range("Q5").Value = ''''''''
$ws.Range("Q5").Style = "Normal"
$ws.Range("R5").Value = ''''''''
$ws.Range("R5").Style = "Normal"
$ws.Range("S5").Value = ''''
$ws.Range("S5").Style = "Normal"

$ws.Range("A6").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = '''4'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = '''CruzHack 2024 - Sitegeist'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''Sitegeist, our technical solution for giving you more problems! We wanted a way to easily analyze trending topics and issues faced by those around us.'
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = '''CruzHack_Sitegeist.md'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = '''2024-1-16'
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = '''Shoto Morisaki'
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = '''project/4/project4-demo.gif'
$ws.Range("I6").Style = "Normal"
$ws.Range("J6").Value = '''Project'
$ws.Range("J6").Style = "Normal"
$ws.Range("K6").Value = '''[ RESTfulAPI, React, Typescript, GPTAPI, HuggingFaceModel, ]'
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").Value = '''""'
$ws.Range("L6").Style = "Normal"
$ws.Range("M6").Value = '''""'
$ws.Range("M6").Style = "Normal"
$ws.Range("N6").Value = '''""'
$ws.Range("N6").Style = "Normal"
$ws.Range("O6").Value = '''""'
$ws.Range("O6").Style = "Normal"
$ws.Range("P6").Value = '''""'
$ws.Range("P6").Style = "Normal"
$ws.Range("Q6").Value = ''''''''
$ws.Range("Q6").Style = "Normal"
$ws.Range("R6").Value = ''''''''
$ws.Range("R6").Style = "Normal"
$ws.Range("S6").Value = ''''
$ws.Range("S6").Style = "Normal"

$ws.Range("A7").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '''12'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '''Machine Learning Internship at LinkX Japan'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''Machine Learning Internship at LinkX Japan, working as Backend / Full Stack engineer. applying RAG system to improve LLM application.'
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = '''ML_internship_LinkX.md'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = '''2024-4-1'
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = '''Shoto Morisaki'
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = '''/project/7/writing.jpg'
$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Value = '''Internship'
$ws.Range("J7").Style = "Normal"
$ws.Range("K7").Value = '''[ Internship, LLM, Natural Language, RAG ]'
$ws.Range("K7").Style = "Normal"
$ws.Range("L7").Value = '''""'
$ws.Range("L7").Style = "Normal"
$ws.Range("M7").Value = '''""'
$ws.Range("M7").Style = "Normal"
$ws.Range("N7").Value = '''""'
$ws.Range("N7").Style = "Normal"
$ws.Range("O7").Value = '''""'
$ws.Range("O7").Style = "Normal"
$ws.Range("P7").Value = '''""'
$ws.Range("P7").Style = "Normal"
$ws.Range("Q7").Value = ''''''''
$ws.Range("Q7").Style = "Normal"
$ws.Range("R7").Value = ''''''''
$ws.Range("R7").Style = "Normal"
$ws.Range("S7").Value = ''''
$ws.Range("S7").Style = "Normal"

$ws.Range("A8").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '''3'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''Qiita × Fast DOCTOR Health Tech Hackathon'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''Medical project to support medical industry with group of University students from California. We had implementation of a visual verification system as a solution. The goal of this approach is to reduce paperwork and streamline processes.'
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = '''Qiita_FastDoctor.md'
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = '''2023-6-14'
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = '''Shoto Morisaki'
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = '''project/3/projectQiitaFD1.jpeg'
$ws.Range("I8").Style = "Normal"
$ws.Range("J8").Value = '''Project'
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").Value = '''[FastAPI, Flask, GPT API]'
$ws.Range("K8").Style = "Normal"
$ws.Range("L8").Value = '''""'
$ws.Range("L8").Style = "Normal"
$ws.Range("M8").Value = '''""'
$ws.Range("M8").Style = "Normal"
$ws.Range("N8").Value = '''""'
$ws.Range("N8").Style = "Normal"
$ws.Range("O8").Value = '''""'
$ws.Range("O8").Style = "Normal"
$ws.Range("P8").Value = '''""'
$ws.Range("P8").Style = "Normal"
$ws.Range("Q8").Value = ''''''''
$ws.Range("Q8").Style = "Normal"
$ws.Range("R8").Value = ''''''''
$ws.Range("R8").Style = "Normal"
$ws.Range("S8").Value = ''''
$ws.Range("S8").Style = "Normal"

$ws.Range("A9").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '''9'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '''Writing README for team project 101'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''When a new member joins a project, writing a README with an overview of the project and links to necessary documents will make it easier to understand the project and reduce unnecessary communication costs.'
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = '''README_for_team.md'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = '''''2024-02-24'''
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = '''Shoto Morisaki'
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Value = ''''
$ws.Range("I9").Style = "Normal"
$ws.Range("J9").Value = '''Blog'
$ws.Range("J9").Style = "Normal"
$ws.Range("K9").Value = '''[ Coding, README ]'
$ws.Range("K9").Style = "Normal"
$ws.Range("L9").Value = '''""'
$ws.Range("L9").Style = "Normal"
$ws.Range("M9").Value = '''""'
$ws.Range("M9").Style = "Normal"
$ws.Range("N9").Value = '''""'
$ws.Range("N9").Style = "Normal"
$ws.Range("O9").Value = '''""'
$ws.Range("O9").Style = "Normal"
$ws.Range("P9").Value = '''""'
$ws.Range("P9").Style = "Normal"
$ws.Range("Q9").Value = ''''''''
$ws.Range("Q9").Style = "Normal"
$ws.Range("R9").Value = ''''''''
$ws.Range("R9").Style = "Normal"
$ws.Range("S9").Value = '''/post/5/sample5.jpg'
$ws.Range("S9").Style = "Normal"

$ws.Range("A10").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = '''6'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '''Comprehensive Survey of Hallucination Mitigation Techniques in Large Language Models'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''A comprehensive survey of techniques developed to reduce hallucinations in LLMs'
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = '''Hallucination_Mitigation_Technique_RAG.md'
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = '''2024-02-11'
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = '''Shoto Morisaki'
$ws.Range("H10").Style = "Normal"
$ws.Range("I10").Value = '''sample/sample4.jpg'
$ws.Range("I10").Style = "Normal"
$ws.Range("J10").Value = '''Research'
$ws.Range("J10").Style = "Normal"
$ws.Range("K10").Value = '''[Large Language Models, Hallucination Mitigation, Techniques, Survey]'
$ws.Range("K10").Style = "Normal"
$ws.Range("L10").Value = '''""'
$ws.Range("L10").Style = "Normal"
$ws.Range("M10").Value = '''""'
$ws.Range("M10").Style = "Normal"
$ws.Range("N10").Value = '''""'
$ws.Range("N10").Style = "Normal"
$ws.Range("O10").Value = '''""'
$ws.Range("O10").Style = "Normal"
$ws.Range("P10").Value = '''""'
$ws.Range("P10").Style = "Normal"
$ws.Range("Q10").Value = ''''''''
$ws.Range("Q10").Style = "Normal"
$ws.Range("R10").Value = ''''''''
$ws.Range("R10").Style = "Normal"
$ws.Range("S10").Value = ''''
$ws.Range("S10").Style = "Normal"

$ws.Range("A11").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = '''14'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '''5 things you shouldn''t do at hackathon and startup'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''When a new member joins a project, writing a README with an overview of the project and links to necessary documents will make it easier to understand the project and reduce unnecessary communication costs.'
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = '''5_things_hackathon_and_startup.md'
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = '''2024-04-8'
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = '''Shoto Morisaki'
$ws.Range("H11").Style = "Normal"
$ws.Range("I11").Value = '''/post/5/sample5.jpg'
$ws.Range("I11").Style = "Normal"
$ws.Range("J11").Value = '''Blog'
$ws.Range("J11").Style = "Normal"
$ws.Range("K11").Value = '''[ Coding, README ]'
$ws.Range("K11").Style = "Normal"
$ws.Range("L11").Value = '''""'
$ws.Range("L11").Style = "Normal"
$ws.Range("M11").Value = '''""'
$ws.Range("M11").Style = "Normal"
$ws.Range("N11").Value = '''""'
$ws.Range("N11").Style = "Normal"
$ws.Range("O11").Value = '''""'
$ws.Range("O11").Style = "Normal"
$ws.Range("P11").Value = '''""'
$ws.Range("P11").Style = "Normal"
$ws.Range("Q11").Value = ''''''''
$ws.Range("Q11").Style = "Normal"
$ws.Range("R11").Value = ''''''''
$ws.Range("R11").Style = "Normal"
$ws.Range("S11").Value = ''''
$ws.Range("S11").Style = "Normal"

$ws.Range("A12").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = '''13'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = '''LLM Law Hackathon Impressions Exploring New Frontiers in Legal Tech'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''Recently, I had the opportunity to participate in the LLM Law Hackathon, an event that brought together legal minds, tech enthusiasts, and innovators to collaborate on creating solutions at the intersection of law and technology. As a participant, I was not only exposed to cutting-edge tools and technologies but also immersed in a vibrant community of like-minded individuals. Here are some of my key impressions from the event.'
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = '''LLm_Law_Hackathon_Stanford.md'
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = '''2024-4-8'
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = '''Shoto Morisaki'
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").Value = '''/project/8/opening.png'
$ws.Range("I12").Style = "Normal"
$ws.Range("J12").Value = '''Hackathon'
$ws.Range("J12").Style = "Normal"
$ws.Range("K12").Value = '''[ Hackathon, LLM, Law, RAG ]'
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = '''""'
$ws.Range("L12").Style = "Normal"
$ws.Range("M12").Value = '''""'
$ws.Range("M12").Style = "Normal"
$ws.Range("N12").Value = '''""'
$ws.Range("N12").Style = "Normal"
$ws.Range("O12").Value = '''""'
$ws.Range("O12").Style = "Normal"
$ws.Range("P12").Value = '''""'
$ws.Range("P12").Style = "Normal"
$ws.Range("Q12").Value = ''''''''
$ws.Range("Q12").Style = "Normal"
$ws.Range("R12").Value = ''''''''
$ws.Range("R12").Style = "Normal"
$ws.Range("S12").Value = ''''
$ws.Range("S12").Style = "Normal"

$ws.Range("A13").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = '''2'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '''Electric Sheep Hackathon (Mobility and Big Data)'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''Collaborated with Japanese university students to create an AR app with Flutter, sponsored by Toyota. Utilized the Luma API and Google Maps API for 3D models and mapping features.'
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = '''Electric_Sheep.md'
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = '''2023-8-14'
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value = '''Shoto Morisaki'
$ws.Range("H13").Style = "Normal"
$ws.Range("I13").Value = '''project/2/projectNagoya2.png'
$ws.Range("I13").Style = "Normal"
$ws.Range("J13").Value = '''Project'
$ws.Range("J13").Style = "Normal"
$ws.Range("K13").Value = '''[ Flutter, 3DModel, GoogleMapAPI, LumaAPI]'
$ws.Range("K13").Style = "Normal"
$ws.Range("L13").Value = '''""'
$ws.Range("L13").Style = "Normal"
$ws.Range("M13").Value = '''""'
$ws.Range("M13").Style = "Normal"
$ws.Range("N13").Value = '''""'
$ws.Range("N13").Style = "Normal"
$ws.Range("O13").Value = '''""'
$ws.Range("O13").Style = "Normal"
$ws.Range("P13").Value = '''""'
$ws.Range("P13").Style = "Normal"
$ws.Range("Q13").Value = ''''''''
$ws.Range("Q13").Style = "Normal"
$ws.Range("R13").Value = ''''''''
$ws.Range("R13").Style = "Normal"
$ws.Range("S13").Value = ''''
$ws.Range("S13").Style = "Normal"

$ws.Range("A14").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '''5'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''Definition of RAG'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''Exploring the RAG Paradigm in LLMs. Technically, RAG uses a variety of innovative approaches that address important questions such as “what to search for,” “when to search for,” and “how to use the information obtained.” It has been strengthened.'
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = '''Definition_of_RAG.md'
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = '''2024-02-3'
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value = '''Shoto Morisaki'
$ws.Range("H14").Style = "Normal"
$ws.Range("I14").Value = '''sample/sample3.jpg'
$ws.Range("I14").Style = "Normal"
$ws.Range("J14").Value = '''Research'
$ws.Range("J14").Style = "Normal"
$ws.Range("K14").Value = '''[RAG, LLM, Knowledge Injection, Fine-tuning, External Data]'
$ws.Range("K14").Style = "Normal"
$ws.Range("L14").Value = '''""'
$ws.Range("L14").Style = "Normal"
$ws.Range("M14").Value = '''""'
$ws.Range("M14").Style = "Normal"
$ws.Range("N14").Value = '''""'
$ws.Range("N14").Style = "Normal"
$ws.Range("O14").Value = '''""'
$ws.Range("O14").Style = "Normal"
$ws.Range("P14").Value = '''""'
$ws.Range("P14").Style = "Normal"
$ws.Range("Q14").Value = ''''''''
$ws.Range("Q14").Style = "Normal"
$ws.Range("R14").Value = ''''''''
$ws.Range("R14").Style = "Normal"
$ws.Range("S14").Value = ''''
$ws.Range("S14").Style = "Normal"

$ws.Range("A15").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = '''''Project'' # ''Project'', ''Article'''
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '''1'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''Symposium Cyberwarfare in Russia and Ukraine'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''Recently, I worked on a school project about the cyberwarfare between Russia and Ukraine. I conducted research on what is happening and how much damage has occurred, as well as what types of cyber attacks have been used. For my research, I primarily referenced a cyber report from Microsoft, as well as cybersecurity conferences like DEF CON and BlackHat.'
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = '''Symposium_Cyberwarfare.md'
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = '''2023-4-21'
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = '''Shoto Morisaki'
$ws.Range("H15").Style = "Normal"
$ws.Range("I15").Value = '''project/1/projectSymposium1.jpeg'
$ws.Range("I15").Style = "Normal"
$ws.Range("J15").Value = '''Project'
$ws.Range("J15").Style = "Normal"
$ws.Range("K15").Value = '''[ Research, Cybersecurity ]'
$ws.Range("K15").Style = "Normal"
$ws.Range("L15").Value = '''""'
$ws.Range("L15").Style = "Normal"
$ws.Range("M15").Value = '''""'
$ws.Range("M15").Style = "Normal"
$ws.Range("N15").Value = '''""'
$ws.Range("N15").Style = "Normal"
$ws.Range("O15").Value = '''""'
$ws.Range("O15").Style = "Normal"
$ws.Range("P15").Value = '''""'
$ws.Range("P15").Style = "Normal"
$ws.Range("Q15").Value = ''''''''
$ws.Range("Q15").Style = "Normal"
$ws.Range("R15").Value = ''''''''
$ws.Range("R15").Style = "Normal"
$ws.Range("S15").Value = ''''
$ws.Range("S15").Style = "Normal"

$ws.Range("A16").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '''8'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''Commonly Used React HooksGithub accounts'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''List of commonly used React hooks with brief explanations.'
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = '''React_Hooks.md'
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = '''2024-02-11'
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = '''Shoto Morisaki'
$ws.Range("H16").Style = "Normal"
$ws.Range("I16").Value = '''sample/opened-laptop.jpg'
$ws.Range("I16").Style = "Normal"
$ws.Range("J16").Value = '''Coding'
$ws.Range("J16").Style = "Normal"
$ws.Range("K16").Value = '''[Git]'
$ws.Range("K16").Style = "Normal"
$ws.Range("L16").Value = '''""'
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = '''""'
$ws.Range("M16").Style = "Normal"
$ws.Range("N16").Value = '''""'
$ws.Range("N16").Style = "Normal"
$ws.Range("O16").Value = '''""'
$ws.Range("O16").Style = "Normal"
$ws.Range("P16").Value = '''""'
$ws.Range("P16").Style = "Normal"
$ws.Range("Q16").Value = ''''''''
$ws.Range("Q16").Style = "Normal"
$ws.Range("R16").Value = ''''''''
$ws.Range("R16").Style = "Normal"
$ws.Range("S16").Value = ''''
$ws.Range("S16").Style = "Normal"

$ws.Range("A17").Value = '''''Published'' # ''Draft'', ''Pending'', ''Published'', ''Rewriting'''
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = '''''Article'' # ''Project'', ''Article'''
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '''7'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''【Note】Procedures for using different ssh connection settings (config) with multiple Github accounts'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''When working with multiple GitHub accounts, it''s crucial to configure SSH connection settings properly to ensure seamless authentication. Here are the steps to set up and manage multiple SSH configurations'
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = '''Different_ssh_settings.md'
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = '''2024-02-11'
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value = '''Shoto Morisaki'
$ws.Range("H17").Style = "Normal"
$ws.Range("I17").Value = '''sample/opened-laptop.jpg'
$ws.Range("I17").Style = "Normal"
$ws.Range("J17").Value = '''Coding'
$ws.Range("J17").Style = "Normal"
$ws.Range("K17").Value = '''[Git]'
$ws.Range("K17").Style = "Normal"
$ws.Range("L17").Value = '''""'
$ws.Range("L17").Style = "Normal"
$ws.Range("M17").Value = '''""'
$ws.Range("M17").Style = "Normal"
$ws.Range("N17").Value = '''""'
$ws.Range("N17").Style = "Normal"
$ws.Range("O17").Value = '''""'
$ws.Range("O17").Style = "Normal"
$ws.Range("P17").Value = '''""'
$ws.Range("P17").Style = "Normal"
$ws.Range("Q17").Value = ''''''''
$ws.Range("Q17").Style = "Normal"
$ws.Range("R17").Value = ''''''''
$ws.Range("R17").Style = "Normal"
$ws.Range("S17").Value = ''''
$ws.Range("S17").Style = "Normal"

